# Sync the "BannedPaths" rule row with its non-localizable counterpart:
# remove the old "BannedPaths" row and add a new "BannedPath" row further
# down the table (with Severity changed from Blocker to Critical and no Tags).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "BannedPaths" rule row (row 35). This shifts rows 36-40 up
# by one (to 35-39).
$ws.Rows(35).Delete()

# Insert a new row at position 40 (after the rows that shifted up) and fill
# it in with the renamed rule.
$ws.Rows(40).Insert()
$ws.Range("A40").Value2 = "BannedPath"
$ws.Range("B40").Value2 = "Customer packages should not install content under /libs"
$ws.Range("C40").Value2 = "Bug"
$ws.Range("D40").Value2 = "Critical"

# Restore the cursor/selection that Excel recorded after the edit.
$ws.Range("A37").Select()
